$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 693.875
$ws.Range("I15").Value = 693.875
$ws.Range("K15").Value = 2081.625
$ws.Range("M15").Value = -1912.625

$ws.Range("H61").Value = 489.8
$ws.Range("I61").Value = 489.8
$ws.Range("K61").Value = 1469.4
$ws.Range("M61").Value = -1297.4

$ws.Range("I107").Value = 20834994
$ws.Range("J107").Value = 992.5
$ws.Range("K107").Value = 20834994
$ws.Range("L107").Value = 992.5
$ws.Range("M107").Value = -20833074
$ws.Range("N107").Value = -4832.5

$ws.Range("H116").Value = 5421.8
$ws.Range("I116").Value = 4370.75
$ws.Range("J116").Value = 6122.5
$ws.Range("K116").Value = 4370.75
$ws.Range("L116").Value = 6122.5
$ws.Range("M116").Value = -928.75
$ws.Range("N116").Value = -13006.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4522.8555
$ws.Range("I32").Value = 2960.6094
$ws.Range("K32").Value = 2960.6094
$ws.Range("M32").Value = -2673.6094

$ws.Range("H52").Value = 99999
$ws.Range("J52").Value = 99999
$ws.Range("L52").Value = 99999
$ws.Range("N52").Value = -100635

$ws.Range("H61").Value = 1758.5
$ws.Range("J61").Value = 2214.4
$ws.Range("L61").Value = 2214.4
$ws.Range("N61").Value = -2638.4

$ws.Range("H132").Value = 2442.3333
$ws.Range("J132").Value = 4163
$ws.Range("L132").Value = 12489
$ws.Range("N132").Value = -17549

$ws.Range("H136").Value = 1758.5
$ws.Range("J136").Value = 2214.4
$ws.Range("L136").Value = 6643.200000000001
$ws.Range("N136").Value = -11743.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2000
$ws.Range("I5").Value = 2000
$ws.Range("K5").Value = 2000
$ws.Range("M5").Value = -1887

$ws.Range("H22").Value = 297.8
$ws.Range("I22").Value = 297.8
$ws.Range("K22").Value = 297.8
$ws.Range("M22").Value = -124.8

$ws.Range("H134").Value = 2830.087
$ws.Range("I134").Value = 1101.2941
$ws.Range("J134").Value = 7728.3335
$ws.Range("K134").Value = 3303.8823
$ws.Range("L134").Value = 23185.0005
$ws.Range("M134").Value = -768.8823000000002
$ws.Range("N134").Value = -28255.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1844.1
$ws.Range("I16").Value = 1430.75
$ws.Range("K16").Value = 1430.75
$ws.Range("M16").Value = -1143.75

$ws.Range("H22").Value = 238.65218
$ws.Range("I22").Value = 163.5
$ws.Range("J22").Value = 739.6667
$ws.Range("K22").Value = 163.5
$ws.Range("L22").Value = 739.6667
$ws.Range("M22").Value = 186.5
$ws.Range("N22").Value = -1439.6667

$ws.Range("H31").Value = 29451.727
$ws.Range("I31").Value = 1276.25
$ws.Range("J31").Value = 72798.62
$ws.Range("K31").Value = 1276.25
$ws.Range("L31").Value = 72798.62
$ws.Range("M31").Value = -981.25
$ws.Range("N31").Value = -73388.62

$ws.Range("H34").Value = 29451.727
$ws.Range("I34").Value = 1276.25
$ws.Range("J34").Value = 72798.62
$ws.Range("K34").Value = 1276.25
$ws.Range("L34").Value = 72798.62
$ws.Range("M34").Value = -1074.25
$ws.Range("N34").Value = -73202.62

$ws.Range("H107").Value = 1771.95
$ws.Range("I107").Value = 1640.3334
$ws.Range("K107").Value = 1640.3334
$ws.Range("M107").Value = 279.6666

$ws.Range("H113").Value = 1844.1
$ws.Range("I113").Value = 1430.75
$ws.Range("K113").Value = 1430.75
$ws.Range("M113").Value = 739.25

$ws.Range("H122").Value = 2657.7693
$ws.Range("I122").Value = 2377.7144
$ws.Range("J122").Value = 2984.5
$ws.Range("K122").Value = 7133.1432
$ws.Range("L122").Value = 8953.5
$ws.Range("M122").Value = -4683.1432
$ws.Range("N122").Value = -13853.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 261.5
$ws.Range("I22").Value = 261.5
$ws.Range("K22").Value = 784.5
$ws.Range("M22").Value = -615.5

$ws.Range("H27").Value = 261.5
$ws.Range("I27").Value = 261.5
$ws.Range("K27").Value = 784.5
$ws.Range("M27").Value = -682.5

$ws.Range("H54").Value = 3000
$ws.Range("J54").Value = 3000
$ws.Range("L54").Value = 9000
$ws.Range("N54").Value = -10118

$ws.Range("H132").Value = 2743.5454
$ws.Range("J132").Value = 3125.7144
$ws.Range("L132").Value = 28131.4296
$ws.Range("N132").Value = -33191.4296

$ws.Range("H137").Value = 3921
$ws.Range("I137").Value = 1405.2
$ws.Range("J137").Value = 5718
$ws.Range("K137").Value = 4215.6
$ws.Range("L137").Value = 17154
$ws.Range("M137").Value = 884.3999999999996
$ws.Range("N137").Value = -27354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 598
$ws.Range("I2").Value = 802.2308
$ws.Range("J2").Value = 67
$ws.Range("K2").Value = 802.2308
$ws.Range("L2").Value = 67
$ws.Range("M2").Value = -689.2308
$ws.Range("N2").Value = -293

$ws.Range("H3").Value = 497.5
$ws.Range("I3").Value = 497.5
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 497.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -381.5
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1706.75
$ws.Range("I16").Value = 1546
$ws.Range("K16").Value = 1546
$ws.Range("M16").Value = -1376

$ws.Range("H61").Value = 6949307.5
$ws.Range("I61").Value = 13895344
$ws.Range("J61").Value = 3271.75
$ws.Range("K61").Value = 13895344
$ws.Range("L61").Value = 3271.75
$ws.Range("M61").Value = -13895142
$ws.Range("N61").Value = -3675.75

$ws.Range("H113").Value = 6949307.5
$ws.Range("I113").Value = 13895344
$ws.Range("J113").Value = 3271.75
$ws.Range("K113").Value = 13895344
$ws.Range("L113").Value = 3271.75
$ws.Range("M113").Value = -13893174
$ws.Range("N113").Value = -7611.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 20000
$ws.Range("J4").Value = 20000
$ws.Range("L4").Value = 20000
$ws.Range("N4").Value = -20226

$ws.Range("H128").Value = 44000
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H129").Value = 39350
$ws.Range("J129").Value = 39700
$ws.Range("L129").Value = 39700
$ws.Range("N129").Value = -49700

